# Generate Report for Handback
#
# The de-de handback has completed (in sync with en-US). This records the
# target/handback files + handback timestamp for both locales, links the
# newly-known target files back to their source docs (like column A
# already does), and widens a few columns so the longer strings fit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$urlFor87e = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cc25fc090a5e4455bf5f9b6b51a5da56247956c/e2e/87e26f2f-6832-4109-8dad-4940bb52adef.md"
$urlForBdb = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cc25fc090a5e4455bf5f9b6b51a5da56247956c/e2e/bdbeb8fd-bd65-4b25-b8e5-eaa7a2691331.md"

$nameFor87e = "87e26f2f-6832-4109-8dad-4940bb52adef.md"
$nameForBdb = "bdbeb8fd-bd65-4b25-b8e5-eaa7a2691331.md"

$hyperlinkColor = 15570276  # BGR for font color FF6495ED used by the existing HyperLink style

# ---------------------------------------------------------------------
# Status text: both locale rows are now handed back in sync with en-US.
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: "Latest Target File" (I) / "Latest Handback File" (J)
# ---------------------------------------------------------------------
$ws2.Range("I2").Value = $nameFor87e
$ws2.Range("J2").Value = "87e26f2f-6832-4109-8dad-4940bb52adef.0924a395015f7fa44de2e12a943391c1c86eb809.zh-cn.xlf"

$ws2.Range("I3").Value = $nameForBdb
$ws2.Range("J3").Value = "bdbeb8fd-bd65-4b25-b8e5-eaa7a2691331.732ce2bbe1788ca57a8948d6caf9063067e9050a.zh-cn.xlf"

$ws2.Hyperlinks.Add($ws2.Range("I2"), $urlFor87e, "", "", $nameFor87e)
$ws2.Hyperlinks.Add($ws2.Range("I3"), $urlForBdb, "", "", $nameForBdb)

$ws2.Range("I2").Font.Underline = $true
$ws2.Range("I2").Font.Color = $hyperlinkColor
$ws2.Range("I3").Font.Underline = $true
$ws2.Range("I3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------
# de-de sheet: "Latest Target File" (I) / "Latest Handback File" (J) /
# "Latest Handback DateTime" (K) -- this is the handback being reported.
# ---------------------------------------------------------------------
$ws3.Range("I2").Value = $nameFor87e
$ws3.Range("J2").Value = "87e26f2f-6832-4109-8dad-4940bb52adef.0924a395015f7fa44de2e12a943391c1c86eb809.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-27 12:58:13"

$ws3.Range("I3").Value = $nameForBdb
$ws3.Range("J3").Value = "bdbeb8fd-bd65-4b25-b8e5-eaa7a2691331.732ce2bbe1788ca57a8948d6caf9063067e9050a.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-27 12:58:13"

$ws3.Hyperlinks.Add($ws3.Range("I2"), $urlFor87e, "", "", $nameFor87e)
$ws3.Hyperlinks.Add($ws3.Range("I3"), $urlForBdb, "", "", $nameForBdb)

$ws3.Range("I2").Font.Underline = $true
$ws3.Range("I2").Font.Color = $hyperlinkColor
$ws3.Range("I3").Font.Underline = $true
$ws3.Range("I3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------
# zh-cn "Latest HO Xliff Generate Date" (column K) records when the
# *handoff* xliff for the zh-cn leg was produced.
# ---------------------------------------------------------------------
$ws2.Range("K2").Value = "2016-08-27 12:58:06"
$ws2.Range("K3").Value = "2016-08-27 12:58:06"

# ---------------------------------------------------------------------
# Column widths: widen "Status"/"Latest Target File"/"Latest Handback
# File" columns so the new longer text is readable.
# ---------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 29.15
$ws1.Columns.Item(6).ColumnWidth = 29.15

$ws2.Columns.Item(3).ColumnWidth = 29.15
$ws2.Columns.Item(9).ColumnWidth = 39.17
$ws2.Columns.Item(10).ColumnWidth = 39.17

$ws3.Columns.Item(3).ColumnWidth = 29.15
$ws3.Columns.Item(9).ColumnWidth = 39.17
$ws3.Columns.Item(10).ColumnWidth = 39.17

Write-Host "applied handback report updates"
